$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the 'author' column (E) for rows 2-9, which previously held "NA"
# placeholders, with the actual scientificNameAuthorship values (matching
# what's already present in column AG for each row).
$ws.Range("E2").Value = "Sutemin, 1969"
$ws.Range("E3").Value = "Sutemin, 1969"
$ws.Range("E4").Value = "Curt., [1829]"
$ws.Range("E5").Value = "Suteminn, 1969"
$ws.Range("E6").Value = "Sutemin, 1969"
$ws.Range("E7").Value = "Kirby, 1837"
$ws.Range("E8").Value = "Westwood, 1858"
$ws.Range("E9").Value = "Alpinus, 1874"

# Match the reviewer's final on-screen selection/scroll position.
$ws.Range("AL1").Select()
